$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 688
$ws1.Range("F3").Value = 53
$ws1.Range("F4").Value = 1998
$ws1.Range("F5").Value = 5843
$ws1.Range("F6").Value = 1648
$ws1.Range("F7").Value = 175
$ws1.Range("F8").Value = 3299
$ws1.Range("F11").Value = 1385
$ws1.Range("F12").Value = 4602
$ws1.Range("F13").Value = 1099
$ws1.Range("F14").Value = 1730
$ws1.Range("F17").Value = 0
$ws1.Range("F18").Value = 58
$ws1.Range("F19").Value = 192
$ws1.Range("F24").Value = 21
$ws1.Range("F29").Value = 1128
$ws1.Range("F30").Value = 419
$ws1.Range("F32").Value = 212
$ws1.Range("F33").Value = 417
$ws1.Range("F34").Value = 994
$ws1.Range("F36").Value = 1766
$ws1.Range("F37").Value = 2272
$ws1.Range("F40").Value = 282
$ws1.Range("F42").Value = 395
$ws1.Range("F43").Value = 44
$ws1.Range("F44").Value = 681
$ws1.Range("F46").Value = 451
$ws1.Range("F47").Value = 427

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 688
$ws4.Range("F4").Value = 53
$ws4.Range("F5").Value = 1998
$ws4.Range("F6").Value = 5843
$ws4.Range("F7").Value = 1648
$ws4.Range("F8").Value = 175
$ws4.Range("F9").Value = 3299
$ws4.Range("F11").Value = 1385
$ws4.Range("F12").Value = 4602
$ws4.Range("F13").Value = 1730
$ws4.Range("F16").Value = 53
$ws4.Range("F19").Value = 58
$ws4.Range("F20").Value = 192
$ws4.Range("F28").Value = 1129
$ws4.Range("F29").Value = 419
$ws4.Range("F31").Value = 212
$ws4.Range("F32").Value = 994
$ws4.Range("F33").Value = 1767
$ws4.Range("F34").Value = 2272
$ws4.Range("F39").Value = 282
$ws4.Range("F41").Value = 395
$ws4.Range("F42").Value = 681
$ws4.Range("F43").Value = 451
$ws4.Range("F44").Value = 427
